# Edit script for bhajans.pptx
#
# Summary of changes (per commit "Fixed line length to 45 instead of 49,
# no overflows now."):
#   1. Remove the last two slides from the deck (slide 3 and slide 4).
#   2. Slide 1: update the date "October 15, 2016" -> "October 22, 2016".
#   3. Slide 2:
#       a. Title textbox: "Aao Aao Aao Aao " -> "Aao Aao Antaryami"
#       b. Body textbox: replace the lyrics with the new (shorter-lined)
#          lyrics text.
#       c. Remove the now-unused trailing placeholder textboxes
#          (TextBox 4, TextBox 5 "Continued", TextBox 6).

$p = $ppt.ActivePresentation

# --- Slide 1: fix the date -------------------------------------------------
$s1 = $p.Slides.Item(1)
# "TextBox 2" holds two runs: the big title, then the smaller date run.
# Editing the individual Run's Text (rather than the whole TextRange)
# keeps each run's own formatting (size/color) and paragraph properties
# intact.
$s1.Shapes.Item(2).TextFrame.TextRange.Runs(2).Text = "October 22, 2016"

# --- Slide 2: update the song title/lyrics ---------------------------------
$s2 = $p.Slides.Item(2)

# "TextBox 2" (the song title line)
$s2.Shapes.Item(2).TextFrame.TextRange.Runs(1).Text = "Aao Aao Antaryami"

# "TextBox 3" (the lyrics body) - single run, multi-line text.
$nl = [char]10
$line1 = "Aao Aao Antaryami Adi Narayana Sai Narayana  "
$line2 = "|Aao Aao...| Hrudaya Vihari Hey Giridhari"
$line3 = "Hey Giridhari Hey Giridhari"
$line4 = "Hari Hari Sri Hari Sai Murari"
$newLyrics = $line1 + $nl + $line2 + $nl + $line3 + $nl + $line4
$s2.Shapes.Item(3).TextFrame.TextRange.Runs(1).Text = $newLyrics

# Remove the trailing placeholder textboxes (TextBox 4, TextBox 5
# "Continued", TextBox 6) - delete from the end so indices stay valid.
$s2.Shapes.Item(6).Delete()
$s2.Shapes.Item(5).Delete()
$s2.Shapes.Item(4).Delete()

# --- Remove the trailing two slides -----------------------------------------
# Delete from the end so the indices of the remaining slides don't shift
# while we're still removing them.
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
